# 2025 roswell data added
# Append two new drug names (glyburide-metformin, glipizide-metformin)
# to the bottom of the SulfonylureaDrugNames list, formatted like the
# console/terminal-pasted data (small monospace font).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new rows at the bottom of the list ---
$ws.Range("A41").Value = "glyburide-metformin"
$ws.Range("A42").Value = "glipizide-metformin"

# --- Style the newly-added cells (pasted-in look: small monospace font) ---
$newRows = $ws.Range("A41:A42")
$newRows.Font.Name = "Lucida Console"
$newRows.Font.Size = 7
$newRows.Font.Color = 0
$newRows.VerticalAlignment = -4108

# --- Restore view state (scrolled down to the newly added rows) ---
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C31").Select()
